$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions): refresh the "想去人数" (want-to-go count) figures
# for a number of still-listed events.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 12591
$wsExpo.Range("F4").Value = 73
$wsExpo.Range("F5").Value = 30
$wsExpo.Range("F7").Value = 10
$wsExpo.Range("F8").Value = 12491
$wsExpo.Range("F9").Value = 249
$wsExpo.Range("F10").Value = 4911
$wsExpo.Range("F11").Value = 4836
$wsExpo.Range("F12").Value = 158
$wsExpo.Range("F15").Value = 112
$wsExpo.Range("F16").Value = 967
$wsExpo.Range("F17").Value = 12
$wsExpo.Range("F20").Value = 79

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances): the 2024-08-25 concert event has ended/dropped
# off the list. Pull the next event's data up into row 2 (keeping the running
# index in column A untouched) and drop the now-duplicate trailing row.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("B3:I3").Copy($wsShow.Range("B2:I2"))
$wsShow.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types): this sheet is the concatenation of
# 演出 + 展览 + 本地生活, so mirror the same two kinds of edits:
#   1) drop the expired 演出 row (row 2) the same way as above
#   2) refresh the same 想去人数 counts for the 展览 rows, shifted down by one
#      row because the 演出 section occupies row 2 in this sheet.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("B3:I23").Copy($wsAll.Range("B2:I22"))
$wsAll.Rows.Item(23).Delete()

$wsAll.Range("F4").Value = 12591
$wsAll.Range("F5").Value = 73
$wsAll.Range("F6").Value = 30
$wsAll.Range("F8").Value = 10
$wsAll.Range("F9").Value = 12491
$wsAll.Range("F10").Value = 249
$wsAll.Range("F11").Value = 4911
$wsAll.Range("F12").Value = 4836
$wsAll.Range("F13").Value = 158
$wsAll.Range("F16").Value = 112
$wsAll.Range("F17").Value = 967
$wsAll.Range("F18").Value = 12
$wsAll.Range("F21").Value = 79
